$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 30 - this shifts the existing rows 30-67 down to 31-68,
# carrying all of their data (and the D-column date style) along with them.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly price-report entry.
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C30").Value = "Arica y Parinacota"
$ws.Range("D30").Value = 44601
$ws.Range("E30").Value = 15
$ws.Range("F30").Value = 100112040
$ws.Range("G30").Value = "Cilantro"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 300
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 1500
$ws.Range("M30").Value = 1250
$ws.Range("N30").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 625
$ws.Range("Q30").Value = 2
$ws.Range("R30").Value = "Hortaliza"
